## Add new registration tests, build and run all tests, fix a spelling error.
## -----------------------------------------------------------------------
## 1) sheet1 (LoginPageModel): select column D (leaves behind the
##    "D1:D1048576" selection seen in the diff) - it will stop being the
##    active tab once the new sheet is activated below.
## 2) Add a new sheet "RegistrationPageModel" after the existing sheet and
##    populate it with the new registration test data, including the two
##    mailto hyperlinks.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(4).Select()

$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "RegistrationPageModel"

# Text-format the RememberMe / ExpectedError columns before typing into them
# so "true" is stored literally and long sentences don't get reinterpreted.
$ws2.Columns.Item(6).NumberFormat = "@"
$ws2.Columns.Item(7).NumberFormat = "@"

# Column widths (best-fit-like, matching the LoginPageModel sheet look).
$ws2.Columns.Item(1).ColumnWidth = 29.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 14.666666666666666
$ws2.Columns.Item(4).ColumnWidth = 9.166666666666666
$ws2.Columns.Item(5).ColumnWidth = 15.833333333333334
$ws2.Columns.Item(6).ColumnWidth = 13.0
$ws2.Columns.Item(7).ColumnWidth = 28.0

# Header row.
$ws2.Range("A1").Value = "Key"
$ws2.Range("B1").Value = "Email"
$ws2.Range("C1").Value = "FullName"
$ws2.Range("D1").Value = "Password"
$ws2.Range("E1").Value = "ConfirmPassword"
$ws2.Range("F1").Value = "RememberMe"
$ws2.Range("G1").Value = "ExpectedError"

# Row 2: NavigateToRegistrationPage - just needs a valid email/password/flag.
$ws2.Range("A2").Value = "NavigateToRegistrationPage"
$ws2.Range("D2").Value = 123456789
$ws2.Range("F2").Value = "true"

# Row 3: RegisterShouldFailWithoutEmail - email left blank on purpose.
$ws2.Range("C3").Value = "Hristina Petkova"
$ws2.Range("A3").Value = "RegisterShouldFailWithoutEmail"
$ws2.Range("D3").Value = 123456789
$ws2.Range("E3").Value = 123456789
$ws2.Range("F3").Value = "true"
$ws2.Range("G3").Value = "The Email field is required."

# Both rows 2 and 4 need the tester's e-mail; add them (and their
# mailto hyperlinks) together.
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:xrissti@gmail.com", "", "", "xrissti@gmail.com")

# Row 4: RegisterShouldFailWithoutFullName - full name left blank on purpose.
$ws2.Range("A4").Value = "RegisterShouldFailWithoutFullName"
$ws2.Hyperlinks.Add($ws2.Range("B4"), "mailto:xrissti@gmail.com", "", "", "xrissti@gmail.com")
$ws2.Range("D4").Value = 123456789
$ws2.Range("E4").Value = 123456789
$ws2.Range("F4").Value = "true"
$ws2.Range("G4").Value = "The Full Name field is required."

# Make the new sheet the active / visible tab, cursor left on G4 like in
# the source file.
$ws2.Range("G4").Select()
$ws2.Activate()
